$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = "eman.tantawi@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg"
$ws.Cells.Item(4, 7).Value = "mariam.noureldin@med.asu.edu.eg, aml.awwad@med.asu.edu.eg, Dr.hend.fariid@med.asu.edu.eg, aya.saeed@med.asu.edu.eg"
$ws.Cells.Item(5, 7).Value = "aml.awwad@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg, wessam.atef@med.asu.edu.eg"
$ws.Cells.Item(7, 7).Value = "afnan.fares@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg"
$ws.Cells.Item(15, 7).Value = "marian.samir@med.asu.edu.eg, nourhan.mostafa@med.asu.edu.eg, dr_ahmed_mostafa@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg"
$ws.Cells.Item(17, 7).Value = "Monica.Eshak@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg"
$ws.Cells.Item(18, 7).Value = "maryam.ashraf@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, marina_atef@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg"
$ws.Cells.Item(19, 7).Value = "eman.tantawi@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg"
$ws.Cells.Item(21, 7).Value = "mariam.noureldin@med.asu.edu.eg, aml.awwad@med.asu.edu.eg, Dr.hend.fariid@med.asu.edu.eg, aya.saeed@med.asu.edu.eg"
$ws.Cells.Item(22, 7).Value = "aml.awwad@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg, wessam.atef@med.asu.edu.eg"
$ws.Cells.Item(24, 7).Value = "afnan.fares@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg"
$ws.Cells.Item(32, 7).Value = "marian.samir@med.asu.edu.eg, nourhan.mostafa@med.asu.edu.eg, dr_ahmed_mostafa@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg"
$ws.Cells.Item(34, 7).Value = "Monica.Eshak@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg"
$ws.Cells.Item(35, 7).Value = "Aya_hamed@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg"
$ws.Cells.Item(45, 7).Value = "backup@backdoor.com, Rania.a.youssef@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg, System"
$ws.Cells.Item(48, 7).Value = "Aya.alaa-eldein@med.asu.edu.eg, marian.samir@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg"
$ws.Cells.Item(51, 7).Value = "yasmin.m.senosy@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, eman.samir@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg"
$ws.Cells.Item(52, 7).Value = "abdullah.elagrody@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg"
$ws.Cells.Item(62, 7).Value = "backup@backdoor.com, Rania.a.youssef@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg, System"
$ws.Cells.Item(65, 7).Value = "Aya.alaa-eldein@med.asu.edu.eg, marian.samir@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg"
$ws.Cells.Item(68, 7).Value = "yasmin.m.senosy@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, eman.samir@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg"
$ws.Cells.Item(69, 7).Value = "abdullah.elagrody@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg"
$ws.Cells.Item(70, 7).Value = "eman.tantawi@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg"
$ws.Cells.Item(72, 7).Value = "Safa.hany@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg, wessam.atef@med.asu.edu.eg"
$ws.Cells.Item(73, 7).Value = "mariam.noureldin@med.asu.edu.eg, aml.awwad@med.asu.edu.eg, Dr.hend.fariid@med.asu.edu.eg, aya.saeed@med.asu.edu.eg"
$ws.Cells.Item(76, 7).Value = "Rania.a.youssef@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg"
$ws.Cells.Item(83, 7).Value = "marian.samir@med.asu.edu.eg, Youstina.ibrahim@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg"
$ws.Cells.Item(85, 7).Value = "maryam.ashraf@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg"
$ws.Cells.Item(86, 7).Value = "maryam.ashraf@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, marina_atef@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg"
$ws.Cells.Item(87, 7).Value = "eman.tantawi@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg"
$ws.Cells.Item(89, 7).Value = "mariam.noureldin@med.asu.edu.eg, aml.awwad@med.asu.edu.eg, Dr.hend.fariid@med.asu.edu.eg, aya.saeed@med.asu.edu.eg"
$ws.Cells.Item(98, 7).Value = "afaf.abdallah@med.asu.edu.eg, nourhanhosni@med.asu.edu.eg, Walaa.h.ghanima@med.asu.edu.eg, user@user.com, nourhanmohamed@med.asu.edu.eg"
$ws.Cells.Item(100, 7).Value = "Aya.alaa-eldein@med.asu.edu.eg, nourhanmohamed@med.asu.edu.eg"
$ws.Cells.Item(102, 7).Value = "maryam.ashraf@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg"
$ws.Cells.Item(103, 7).Value = "neveen.nashaat@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg"
$ws.Cells.Item(105, 7).Value = "backup@backdoor.com, System"
$ws.Cells.Item(106, 7).Value = "aml.awwad@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg, wessam.atef@med.asu.edu.eg"
$ws.Cells.Item(109, 7).Value = "afnan.fares@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg"
$ws.Cells.Item(110, 7).Value = "afnan.fares@med.asu.edu.eg, Alshimaa_khaled@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg"
$ws.Cells.Item(116, 7).Value = "nourhan.mostafa@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg, enas.omran@med.asu.edu.eg"
$ws.Cells.Item(117, 7).Value = "user@user.com, takyas@med.asu.edu.eg, Amr-Saeed@med.asu.edu.eg`n, enas.omran@med.asu.edu.eg"
$ws.Cells.Item(119, 7).Value = "yasmin.m.senosy@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, marina_atef@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, marinasorial@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg"
$ws.Cells.Item(120, 7).Value = "neveen.nashaat@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg"
$ws.Cells.Item(122, 7).Value = "backup@backdoor.com, System"
$ws.Cells.Item(123, 7).Value = "aml.awwad@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg, wessam.atef@med.asu.edu.eg"
$ws.Cells.Item(126, 7).Value = "afnan.fares@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg"
$ws.Cells.Item(127, 7).Value = "afnan.fares@med.asu.edu.eg, Alshimaa_khaled@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg"
$ws.Cells.Item(133, 7).Value = "nourhan.mostafa@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg, enas.omran@med.asu.edu.eg"
$ws.Cells.Item(134, 7).Value = "user@user.com, takyas@med.asu.edu.eg, Amr-Saeed@med.asu.edu.eg`n, enas.omran@med.asu.edu.eg"
$ws.Cells.Item(136, 7).Value = "yasmin.m.senosy@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, marina_atef@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, marinasorial@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg"
$ws.Cells.Item(137, 7).Value = "neveen.nashaat@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg"
$ws.Cells.Item(143, 7).Value = "afnan.fares@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg"
$ws.Cells.Item(144, 7).Value = "nourhanmamdouh@med.asu.edu.eg, Mayarembaby@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg"
$ws.Cells.Item(150, 7).Value = "marian.samir@med.asu.edu.eg, Youstina.ibrahim@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg"
$ws.Cells.Item(151, 7).Value = "marian.samir@med.asu.edu.eg, nourhanmohamed@med.asu.edu.eg, randa.rabea@med.asu.edu.eg, System"
